# Word COM-interop script implementing:
#  1. Remove the standalone "Meta description: ..." paragraph that used to
#     follow the title heading.
#  2. Replace the final paragraph (the DALL-E image-prompt text) with two
#     new paragraphs: a bold repeat of the title, followed by an italic
#     paragraph containing the meta-description copy (minus the
#     "Meta description" label).

$d = $word.ActiveDocument

# --- Step 1: delete the "Meta description" paragraph (paragraph #2) ---
# Locate it robustly via Find rather than assuming a fixed index.
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("Meta description*rounds.", $true, $false, $true, $false, $false, `
                                  $true, 1, $false, "", 0)
if ($found) {
    # Extend the range to also swallow the paragraph mark so the
    # paragraph itself disappears rather than leaving a blank line.
    $paraRange = $d.Paragraphs.Item($findRange.Paragraphs.First.Index).Range
    $paraRange.Delete()
} else {
    # Fallback: the paragraph right after the title heading.
    $d.Paragraphs.Item(2).Range.Delete()
}

# --- Step 2: replace the closing DALL-E prompt paragraph ---
# Find the paragraph that holds the DALL-E image-prompt text; fall back to
# "the last paragraph in the document" if, for some reason, it can't be
# located by content.
$dalleRange = $d.Content
$dalleRange.Find.ClearFormatting()
$dalleFound = $dalleRange.Find.Execute("Create a feature image for*exciting theme and graphics.", `
                                        $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
if ($dalleFound) {
    $targetParaIndex = $dalleRange.Paragraphs.First.Index
    $lastPara = $d.Paragraphs.Item($targetParaIndex)
} else {
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
}
$lastRange = $lastPara.Range

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newXml = "<w:p $ns><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruit Spin for Free - Classic Slot Game with Modern Features</w:t></w:r></w:p>" + `
          "<w:p $ns><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Fruit Spin, a classic slot game with modern features. Play for free and enjoy the chance to win cash prizes and trigger bonus rounds.</w:t></w:r></w:p>"

[void]$lastRange.InsertXML($newXml)
